$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 944.6875
$ws.Range("J97").Value = 943.9286
$ws.Range("L97").Value = 2831.7858
$ws.Range("N97").Value = -3823.7858
$ws.Range("H116").Value = 7632.0454
$ws.Range("I116").Value = 10292.357
$ws.Range("J116").Value = 2976.5
$ws.Range("K116").Value = 10292.357
$ws.Range("L116").Value = 2976.5
$ws.Range("M116").Value = -6850.357
$ws.Range("N116").Value = -9860.5
$ws.Range("H132").Value = 2909.025
$ws.Range("I132").Value = 2774.4866
$ws.Range("K132").Value = 8323.459800000001
$ws.Range("M132").Value = -5793.459800000001
$ws.Range("H137").Value = 1939611.4
$ws.Range("I137").Value = 3206284.8
$ws.Range("J137").Value = 2346
$ws.Range("K137").Value = 9618854.399999999
$ws.Range("L137").Value = 7038
$ws.Range("M137").Value = -9616304.399999999
$ws.Range("N137").Value = -12138
$ws.Range("H138").Value = 4083.4285
$ws.Range("I138").Value = 3492.5
$ws.Range("J138").Value = 4293.113
$ws.Range("K138").Value = 10477.5
$ws.Range("L138").Value = 12879.339
$ws.Range("M138").Value = -5337.5
$ws.Range("N138").Value = -23159.339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10419330
$ws.Range("I61").Value = 16668249
$ws.Range("J61").Value = 4466.6665
$ws.Range("K61").Value = 16668249
$ws.Range("L61").Value = 4466.6665
$ws.Range("M61").Value = -16668037
$ws.Range("N61").Value = -4890.6665
$ws.Range("H74").Value = 11365328
$ws.Range("I74").Value = 990.89655
$ws.Range("J74").Value = 33336382
$ws.Range("K74").Value = 990.89655
$ws.Range("L74").Value = 33336382
$ws.Range("M74").Value = -116.89655
$ws.Range("N74").Value = -33338130
$ws.Range("H77").Value = 11365328
$ws.Range("I77").Value = 990.89655
$ws.Range("J77").Value = 33336382
$ws.Range("K77").Value = 4954.48275
$ws.Range("L77").Value = 166681910
$ws.Range("M77").Value = -586.4827500000001
$ws.Range("N77").Value = -166690646
$ws.Range("H88").Value = 4633.3335
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 6150
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 6150
$ws.Range("M88").Value = -1194
$ws.Range("N88").Value = -6962
$ws.Range("H91").Value = 4633.3335
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 6150
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 6150
$ws.Range("M91").Value = -196
$ws.Range("N91").Value = -8958
$ws.Range("H122").Value = 168666.5
$ws.Range("I122").Value = 201799.8
$ws.Range("K122").Value = 605399.3999999999
$ws.Range("M122").Value = -602949.3999999999
$ws.Range("H132").Value = 1482806
$ws.Range("I132").Value = 2794.818
$ws.Range("J132").Value = 4053351.8
$ws.Range("K132").Value = 8384.454000000002
$ws.Range("L132").Value = 12160055.4
$ws.Range("M132").Value = -5854.454000000002
$ws.Range("N132").Value = -12165115.4
$ws.Range("H136").Value = 10419330
$ws.Range("I136").Value = 16668249
$ws.Range("J136").Value = 4466.6665
$ws.Range("K136").Value = 50004747
$ws.Range("L136").Value = 13399.9995
$ws.Range("M136").Value = -50002197
$ws.Range("N136").Value = -18499.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 275
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -127
$ws.Range("N22").Value = -596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 533.38464
$ws.Range("I22").Value = 378.57144
$ws.Range("J22").Value = 714
$ws.Range("K22").Value = 378.57144
$ws.Range("L22").Value = 714
$ws.Range("M22").Value = -28.57144
$ws.Range("N22").Value = -1414
$ws.Range("H31").Value = 8114.864
$ws.Range("I31").Value = 3654.8948
$ws.Range("K31").Value = 3654.8948
$ws.Range("M31").Value = -3359.8948
$ws.Range("H34").Value = 8114.864
$ws.Range("I34").Value = 3654.8948
$ws.Range("K34").Value = 3654.8948
$ws.Range("M34").Value = -3452.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5670.4165
$ws.Range("I131").Value = 449.7
$ws.Range("J131").Value = 9399.5
$ws.Range("K131").Value = 1349.1
$ws.Range("L131").Value = 28198.5
$ws.Range("M131").Value = 3690.9
$ws.Range("N131").Value = -38278.5
$ws.Range("H141").Value = 6673.7095
$ws.Range("I141").Value = 3133.077
$ws.Range("J141").Value = 9230.833000000001
$ws.Range("K141").Value = 9399.231
$ws.Range("L141").Value = 27692.499
$ws.Range("M141").Value = -4219.231
$ws.Range("N141").Value = -38052.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3526
$ws.Range("I122").Value = 7150
$ws.Range("J122").Value = 2318
$ws.Range("K122").Value = 21450
$ws.Range("L122").Value = 6954
$ws.Range("M122").Value = -19000
$ws.Range("N122").Value = -11854
$ws.Range("H132").Value = 23260116
$ws.Range("I132").Value = 37042244
$ws.Range("J132").Value = 2774.875
$ws.Range("K132").Value = 111126732
$ws.Range("L132").Value = 8324.625
$ws.Range("M132").Value = -111124202
$ws.Range("N132").Value = -13384.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5506.769
$ws.Range("I93").Value = 6789.0527
$ws.Range("J93").Value = 2026.2858
$ws.Range("K93").Value = 6789.0527
$ws.Range("L93").Value = 2026.2858
$ws.Range("M93").Value = -5541.0527
$ws.Range("N93").Value = -4522.2858
$ws.Range("H122").Value = 5205
$ws.Range("I122").Value = 4330.4614
$ws.Range("J122").Value = 5642.269
$ws.Range("K122").Value = 12991.3842
$ws.Range("L122").Value = 16926.807
$ws.Range("M122").Value = -10541.3842
$ws.Range("N122").Value = -21826.807
$ws.Range("H136").Value = 5377831.5
$ws.Range("I136").Value = 1243.64
$ws.Range("J136").Value = 27780280
$ws.Range("K136").Value = 3730.92
$ws.Range("L136").Value = 83340840
$ws.Range("M136").Value = -1180.92
$ws.Range("N136").Value = -83345940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 502500000
$ws.Range("I2").Value = 5000000
$ws.Range("J2").Value = 1000000000
$ws.Range("K2").Value = 5000000
$ws.Range("L2").Value = 1000000000
$ws.Range("M2").Value = -4999888
$ws.Range("N2").Value = -1000000224
$ws.Range("H41").Value = 5550
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 5550
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 5550
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -6330
$ws.Range("H45").Value = 12125.4
$ws.Range("J45").Value = 12125.4
$ws.Range("L45").Value = 12125.4
$ws.Range("N45").Value = -13107.4
$ws.Range("H108").Value = 110626
$ws.Range("J108").Value = 110626
$ws.Range("L108").Value = 110626
$ws.Range("N108").Value = -118306
$ws.Range("H122").Value = 2455.9312
$ws.Range("I122").Value = 2097.1
$ws.Range("J122").Value = 3253.3333
$ws.Range("K122").Value = 6291.299999999999
$ws.Range("L122").Value = 9759.999899999999
$ws.Range("M122").Value = -3841.299999999999
$ws.Range("N122").Value = -14659.9999
